$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2.0
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.7519083333333333
$ws.Range("H2").Value = 2.255725
$ws.Range("I2").Value = 0.135467181536743
$ws.Range("J2").Value = 0.135467181536743
$ws.Range("M2").Value = 2.5191905
$ws.Range("N2").Value = 5.038381
$ws.Range("O2").Value = 0.2895559137428469
$ws.Range("P2").Value = 0.2218380148416813
$ws.Range("Q2").Value = 1.894200330204167
$ws.Range("R2").Value = 11.365201981225
$ws.Range("S2").Value = 0.03922532353203975
$ws.Range("T2").Value = 0.03005177062830874
$ws.Range("E3").Value = 2.0
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.7519083333333333
$ws.Range("H3").Value = 2.255725
$ws.Range("I3").Value = 0.135467181536743
$ws.Range("J3").Value = 0.135467181536743
$ws.Range("O3").Value = 0.1720786450558915
$ws.Range("P3").Value = 0.1977524022342039
$ws.Range("Q3").Value = 1.125694247002778
$ws.Range("R3").Value = 10.131248223025
$ws.Range("S3").Value = 0.02331100904838323
$ws.Range("T3").Value = 0.02678896057278793
$ws.Range("E4").Value = 2.0
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.7519083333333333
$ws.Range("H4").Value = 2.255725
$ws.Range("I4").Value = 0.135467181536743
$ws.Range("J4").Value = 0.135467181536743
$ws.Range("M4").Value = 0.593622
$ws.Range("N4").Value = 1.780866
$ws.Range("O4").Value = 0.06823094983402654
$ws.Range("P4").Value = 0.07841085819810882
$ws.Range("Q4").Value = 0.44634932865
$ws.Range("R4").Value = 4.01714395785
$ws.Range("S4").Value = 0.00924305446759048
$ws.Range("T4").Value = 0.01062209796197502
$ws.Range("E5").Value = 2.0
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.7519083333333333
$ws.Range("H5").Value = 2.255725
$ws.Range("I5").Value = 0.135467181536743
$ws.Range("J5").Value = 0.135467181536743
$ws.Range("M5").Value = 0.8693875
$ws.Range("N5").Value = 1.738775
$ws.Range("O5").Value = 0.09992745366382942
$ws.Range("P5").Value = 0.07655760734576136
$ws.Range("Q5").Value = 0.6536997061458333
$ws.Range("R5").Value = 3.922198236875
$ws.Range("S5").Value = 0.01353689050598246
$ws.Range("T5").Value = 0.01037104329232695
$ws.Range("E6").Value = 2.0
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.7519083333333333
$ws.Range("H6").Value = 2.255725
$ws.Range("I6").Value = 0.135467181536743
$ws.Range("J6").Value = 0.135467181536743
$ws.Range("M6").Value = 2.794049333333334
$ws.Range("N6").Value = 8.382148
$ws.Range("O6").Value = 0.321148205249236
$ws.Range("P6").Value = 0.3690628145090992
$ws.Range("Q6").Value = 2.100868977477778
$ws.Range("R6").Value = 18.9078207973
$ws.Range("S6").Value = 0.04350504222069747
$ws.Range("T6").Value = 0.04999589929156546
$ws.Range("E7").Value = 2.0
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.7519083333333333
$ws.Range("H7").Value = 2.255725
$ws.Range("I7").Value = 0.135467181536743
$ws.Range("J7").Value = 0.135467181536743
$ws.Range("M7").Value = 0.426821
$ws.Range("N7").Value = 1.280463
$ws.Range("O7").Value = 0.04905883245416956
$ws.Range("P7").Value = 0.05637830287114527
$ws.Range("Q7").Value = 0.3209302667416666
$ws.Range("R7").Value = 2.888372400675
$ws.Range("S7").Value = 0.006645861762049649
$ws.Range("T7").Value = 0.007637409789778917
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 0.5089346666666666
$ws.Range("H8").Value = 1.526804
$ws.Range("I8").Value = 0.09169195475469101
$ws.Range("J8").Value = 0.09169195475469102
$ws.Range("M8").Value = 2.5191905
$ws.Range("N8").Value = 5.038381
$ws.Range("O8").Value = 0.2895559137428469
$ws.Range("P8").Value = 0.2218380148416813
$ws.Range("Q8").Value = 1.282103377387333
$ws.Range("R8").Value = 7.692620264324
$ws.Range("S8").Value = 0.02654994774186233
$ws.Range("T8").Value = 0.02034076121973392
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 0.5089346666666666
$ws.Range("H9").Value = 1.526804
$ws.Range("I9").Value = 0.09169195475469101
$ws.Range("J9").Value = 0.09169195475469102
$ws.Range("O9").Value = 0.1720786450558915
$ws.Range("P9").Value = 0.1977524022342039
$ws.Range("Q9").Value = 0.7619344020662222
$ws.Range("R9").Value = 6.857409618595999
$ws.Range("S9").Value = 0.01577822733671334
$ws.Range("T9").Value = 0.01813230431829009
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 0.5089346666666666
$ws.Range("H10").Value = 1.526804
$ws.Range("I10").Value = 0.09169195475469101
$ws.Range("J10").Value = 0.09169195475469102
$ws.Range("M10").Value = 0.593622
$ws.Range("N10").Value = 1.780866
$ws.Range("O10").Value = 0.06823094983402654
$ws.Range("P10").Value = 0.07841085819810882
$ws.Range("Q10").Value = 0.302114814696
$ws.Range("R10").Value = 2.719033332264
$ws.Range("S10").Value = 0.006256229165051154
$ws.Range("T10").Value = 0.007189644862177488
$ws.Range("E11").Value = 3.0
$ws.Range("F11").Value = 1.0
$ws.Range("G11").Value = 0.5089346666666666
$ws.Range("H11").Value = 1.526804
$ws.Range("I11").Value = 0.09169195475469101
$ws.Range("J11").Value = 0.09169195475469102
$ws.Range("M11").Value = 0.8693875
$ws.Range("N11").Value = 1.738775
$ws.Range("O11").Value = 0.09992745366382942
$ws.Range("P11").Value = 0.07655760734576136
$ws.Range("Q11").Value = 0.4424614375166666
$ws.Range("R11").Value = 2.6547686251
$ws.Range("S11").Value = 0.00916254356009533
$ws.Range("T11").Value = 0.007019716668874952
$ws.Range("E12").Value = 3.0
$ws.Range("F12").Value = 1.0
$ws.Range("G12").Value = 0.5089346666666666
$ws.Range("H12").Value = 1.526804
$ws.Range("I12").Value = 0.09169195475469101
$ws.Range("J12").Value = 0.09169195475469102
$ws.Range("M12").Value = 2.794049333333334
$ws.Range("N12").Value = 8.382148
$ws.Range("O12").Value = 0.321148205249236
$ws.Range("P12").Value = 0.3690628145090992
$ws.Range("Q12").Value = 1.421988566110222
$ws.Range("R12").Value = 12.797897094992
$ws.Range("S12").Value = 0.02944670670526317
$ws.Range("T12").Value = 0.03384009088960725
$ws.Range("E13").Value = 3.0
$ws.Range("F13").Value = 1.0
$ws.Range("G13").Value = 0.5089346666666666
$ws.Range("H13").Value = 1.526804
$ws.Range("I13").Value = 0.09169195475469101
$ws.Range("J13").Value = 0.09169195475469102
$ws.Range("M13").Value = 0.426821
$ws.Range("N13").Value = 1.280463
$ws.Range("O13").Value = 0.04905883245416956
$ws.Range("P13").Value = 0.05637830287114527
$ws.Range("Q13").Value = 0.2172240033613333
$ws.Range("R13").Value = 1.955016030252
$ws.Range("S13").Value = 0.004498300245705683
$ws.Range("T13").Value = 0.005169436796007319
$ws.Range("G14").Value = 4.289640333333334
$ws.Range("H14").Value = 12.868921
$ws.Range("I14").Value = 0.7728408637085659
$ws.Range("J14").Value = 0.7728408637085659
$ws.Range("M14").Value = 2.5191905
$ws.Range("N14").Value = 5.038381
$ws.Range("O14").Value = 0.2895559137428469
$ws.Range("P14").Value = 0.2218380148416813
$ws.Range("Q14").Value = 10.80642117615017
$ws.Range("R14").Value = 64.838527056901
$ws.Range("S14").Value = 0.2237806424689448
$ws.Range("T14").Value = 0.1714454829936387
$ws.Range("G15").Value = 4.289640333333334
$ws.Range("H15").Value = 12.868921
$ws.Range("I15").Value = 0.7728408637085659
$ws.Range("J15").Value = 0.7728408637085659
$ws.Range("O15").Value = 0.1720786450558915
$ws.Range("P15").Value = 0.1977524022342039
$ws.Range("Q15").Value = 6.422090607158778
$ws.Range("R15").Value = 57.798815464429
$ws.Range("S15").Value = 0.1329894086707949
$ws.Range("T15").Value = 0.1528311373431259
$ws.Range("G16").Value = 4.289640333333334
$ws.Range("H16").Value = 12.868921
$ws.Range("I16").Value = 0.7728408637085659
$ws.Range("J16").Value = 0.7728408637085659
$ws.Range("M16").Value = 0.593622
$ws.Range("N16").Value = 1.780866
$ws.Range("O16").Value = 0.06823094983402654
$ws.Range("P16").Value = 0.07841085819810882
$ws.Range("Q16").Value = 2.546424873954
$ws.Range("R16").Value = 22.917823865586
$ws.Range("S16").Value = 0.0527316662013849
$ws.Range("T16").Value = 0.06059911537395631
$ws.Range("G17").Value = 4.289640333333334
$ws.Range("H17").Value = 12.868921
$ws.Range("I17").Value = 0.7728408637085659
$ws.Range("J17").Value = 0.7728408637085659
$ws.Range("M17").Value = 0.8693875
$ws.Range("N17").Value = 1.738775
$ws.Range("O17").Value = 0.09992745366382942
$ws.Range("P17").Value = 0.07655760734576136
$ws.Range("Q17").Value = 3.729359685295834
$ws.Range("R17").Value = 22.376158111775
$ws.Range("S17").Value = 0.07722801959775162
$ws.Range("T17").Value = 0.05916684738455946
$ws.Range("G18").Value = 4.289640333333334
$ws.Range("H18").Value = 12.868921
$ws.Range("I18").Value = 0.7728408637085659
$ws.Range("J18").Value = 0.7728408637085659
$ws.Range("M18").Value = 2.794049333333334
$ws.Range("N18").Value = 8.382148
$ws.Range("O18").Value = 0.321148205249236
$ws.Range("P18").Value = 0.3690628145090992
$ws.Range("Q18").Value = 11.98546671358978
$ws.Range("R18").Value = 107.869200422308
$ws.Range("S18").Value = 0.2481964563232754
$ws.Range("T18").Value = 0.2852268243279265
$ws.Range("G19").Value = 4.289640333333334
$ws.Range("H19").Value = 12.868921
$ws.Range("I19").Value = 0.7728408637085659
$ws.Range("J19").Value = 0.7728408637085659
$ws.Range("M19").Value = 0.426821
$ws.Range("N19").Value = 1.280463
$ws.Range("O19").Value = 0.04905883245416956
$ws.Range("P19").Value = 0.05637830287114527
$ws.Range("Q19").Value = 1.830908576713667
$ws.Range("R19").Value = 16.478177190423
$ws.Range("S19").Value = 0.03791467044641422
$ws.Range("T19").Value = 0.04357145628535903
